$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Unbounded KnapSack" pattern, added in column K (mirrors the existing
#     A/C/E/G/I header-then-two-examples layout) ---
$ws.Range("K1").Value = "Unbounded KnapSack"
$ws.Range("K2").Value = "Leetcode - 519"
$ws.Range("K3").Value = "Leetcode - 322"

# Copy the bold header formatting from an existing pattern header onto the new one
$ws.Range("A1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null

# --- New "KeyPad Numeric" pattern, added below the existing table starting at row 7 ---
# (written in shared-string order: the example row first, then the header)
$ws.Range("A8").Value = "LeetCode - 935"
$ws.Range("A7").Value = "KeyPad Numeric"

$ws.Range("A1").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(7).RowHeight = 15.6

# Column widths settle close to their real best-fit values for the new content
$ws.Columns.Item(11).ColumnWidth = 21
$ws.Columns.Item(13).ColumnWidth = 12.666667

# Move the active selection, matching the author's final cursor position
$ws.Range("C13").Select() | Out-Null
